$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.ClearFormats()
}

$ws.Range("D2").Value = "64.512.42"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "3.350.69"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue $ws.Range("D5") "558.45"
$ws.Range("E5").Value = "  -0.17%  "
Set-TextValue $ws.Range("D6") "175.79"
$ws.Range("E6").Value = "  -0.76%  "
Set-TextValue $ws.Range("D7") "0.620"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("D8").Value = "3.341.90"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.05%  "
Set-TextValue $ws.Range("D10") "0.166"
$ws.Range("E10").Value = "  +4.50%  "
Set-TextValue $ws.Range("D11") "0.633"
$ws.Range("E11").Value = "  +1.51%  "
Set-TextValue $ws.Range("D12") "53.83"
$ws.Range("E12").Value = "  +0.14%  "
Set-TextValue $ws.Range("D13") "0.0000275"
$ws.Range("E13").Value = "  +2.29%  "
Set-TextValue $ws.Range("D14") "9.08"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").Value = "3.896.78"
$ws.Range("E15").Value = "  +0.48%  "
Set-TextValue $ws.Range("D16") "18.22"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "3.355.09"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "64.563.07"
$ws.Range("E19").Value = "  +0.92%  "
Set-TextValue $ws.Range("D20") "11.75"
$ws.Range("E20").Value = "  -0.97%  "
Set-TextValue $ws.Range("D21") "0.988"
$ws.Range("E21").Value = "  +0.26%  "
Set-TextValue $ws.Range("D22") "451.40"
$ws.Range("E22").Value = "  +0.78%  "
Set-TextValue $ws.Range("D23") "4.93"
$ws.Range("E23").Value = "  +8.44%  "
Set-TextValue $ws.Range("D24") "4.12"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D25") "14.20"
$ws.Range("E25").Value = "  +7.45%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D26") "86.76"
$ws.Range("E26").Value = "  +2.04%  "
Set-TextValue $ws.Range("D27") "2.90"
$ws.Range("E27").Value = "  +2.27%  "
Set-TextValue $ws.Range("D28") "10.76"
$ws.Range("E28").Value = "  +0.30%  "
Set-TextValue $ws.Range("D29") "8.72"
$ws.Range("E29").Value = "  -0.72%  "
Set-TextValue $ws.Range("D30") "30.83"
$ws.Range("E30").Value = "  +3.80%  "
Set-TextValue $ws.Range("D31") "6.63"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("E32").Value = "  -0.35%  "
Set-TextValue $ws.Range("D33") "571.42"
$ws.Range("E33").Value = "  -3.01%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D34") "60.98"
$ws.Range("E34").Value = "  +4.03%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D35") "0.107"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("B36").Value = "Dai"
$ws.Range("C36").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D36") "1.00"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D37") "3.64"
$ws.Range("E37").Value = "  +3.87%  "
Set-TextValue $ws.Range("D38") "0.140"
$ws.Range("E38").Value = "  -3.65%  "
Set-TextValue $ws.Range("D39") "35.40"
$ws.Range("E39").Value = "  -0.79%  "
Set-TextValue $ws.Range("D40") "0.369"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").Value = "0.0₃0739"
$ws.Range("E41").Value = "  -1.49%  "
$ws.Range("D42").Value = "3.069.74"
$ws.Range("E42").Value = "  -2.03%  "
Set-TextValue $ws.Range("D43") "2.81"
$ws.Range("E43").Value = "  -1.82%  "
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("E46").Value = "  -0.33%  "
Set-TextValue $ws.Range("D47") "3.14"
$ws.Range("E47").Value = "  -1.51%  "
Set-TextValue $ws.Range("D48") "1.00"
$ws.Range("E48").Value = "  +0.19%  "
Set-TextValue $ws.Range("D49") "139.46"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("E51").Value = "  -0.29%  "
